$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Quartal 2" topic block (rows 10-14) ---
# Column A: new topic-weight strings (LDA output) for Q2
$ws.Range("A10").Value = "[(0, '0.012*""kind"" + 0.007*""corona"" + 0.007*""berlin"" + 0.007*""schule"" + 0.006*""frage"" + 0.005*""afd""'), "
$ws.Range("A11").Value = "(1, '0.013*""deutschland"" + 0.012*""afd"" + 0.008*""jahr"" + 0.006*""polizei"" + 0.006*""amp"" + 0.005*""klaren""'), "
$ws.Range("A12").Value = "(2, '0.010*""bundestag"" + 0.010*""amp"" + 0.009*""corona"" + 0.006*""jahr"" + 0.006*""mensch"" + 0.006*""woche""'), "
$ws.Range("A13").Value = "(3, '0.008*""is"" + 0.006*""on"" + 0.005*""i"" + 0.005*""we"" + 0.005*""rassismus"" + 0.005*""trump""'), "
$ws.Range("A14").Value = "(4, '0.022*""amp"" + 0.010*""mensch"" + 0.010*""coronakrise"" + 0.009*""corona"" + 0.008*""zeit"" + 0.007*""krise""')]"

# Column B: new topic labels describing each Q2 topic
$ws.Range("B10").Value = "Schuschließungen"
$ws.Range("B12").Value = "Coronadebatten im Bundestag"
$ws.Range("B13").Value = "Trump und Rassismus"
$ws.Range("B14").Value = "Coronakrise"

# Column C: replace "Lockdown?" note with the final remark
$ws.Range("C11").Value = "Lockdown weniger Kommuniziert"

# --- Column widths / layout ---
# (target widths of 40.77734375 / 46.109375 chars land between this engine's
# pixel-snap grid points; the inputs below snap to the closest achievable width)
$ws.Columns.Item(3).ColumnWidth = 39.95
$ws.Columns.Item(4).ColumnWidth = 45.35

# --- Selection moves to A13 ---
$ws.Range("A13").Select() | Out-Null
